$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(4, 6).Value = 2813
$ws.Cells.Item(6, 6).Value = 222
$ws.Cells.Item(9, 6).Value = 6522
$ws.Cells.Item(12, 6).Value = 79
$ws.Cells.Item(13, 6).Value = 5083
$ws.Cells.Item(14, 6).Value = 4
$ws.Cells.Item(15, 6).Value = 563
$ws.Cells.Item(16, 6).Value = 2679
$ws.Cells.Item(17, 6).Value = 1368
$ws.Cells.Item(18, 6).Value = 1531
$ws.Cells.Item(19, 6).Value = 1241
$ws.Cells.Item(21, 6).Value = 130
$ws.Cells.Item(22, 6).Value = 143
$ws.Cells.Item(23, 6).Value = 1110
$ws.Cells.Item(24, 6).Value = 259
$ws.Cells.Item(25, 6).Value = 552
$ws.Cells.Item(26, 6).Value = 1400
$ws.Cells.Item(28, 6).Value = 2122
$ws.Cells.Item(29, 6).Value = 599
$ws.Cells.Item(30, 6).Value = 49
$ws.Cells.Item(31, 6).Value = 44
$ws.Cells.Item(32, 6).Value = 117
$ws.Cells.Item(33, 6).Value = 267
$ws.Cells.Item(34, 6).Value = 1543
$ws.Cells.Item(38, 6).Value = 1101
$ws.Cells.Item(41, 6).Value = 2329
$ws.Cells.Item(42, 6).Value = 2593
$ws.Cells.Item(43, 6).Value = 58
$ws.Cells.Item(44, 6).Value = 152
$ws.Cells.Item(46, 6).Value = 282
$ws.Cells.Item(48, 6).Value = 115

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(9, 6).Value = 332
$ws.Cells.Item(13, 6).Value = 204
$ws.Cells.Item(16, 6).Value = 243
$ws.Cells.Item(17, 6).Value = 166
$ws.Cells.Item(27, 6).Value = 427
$ws.Cells.Item(28, 6).Value = 31

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(4, 6).Value = 520
$ws.Cells.Item(6, 6).Value = 1711
$ws.Cells.Item(7, 6).Value = 575
$ws.Cells.Item(8, 6).Value = 1550
$ws.Cells.Item(9, 6).Value = 1826
$ws.Cells.Item(10, 6).Value = 2593
$ws.Cells.Item(11, 6).Value = 899
$ws.Cells.Item(12, 6).Value = 785
$ws.Cells.Item(14, 6).Value = 159

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 76
$ws.Cells.Item(3, 6).Value = 520
$ws.Cells.Item(4, 6).Value = 1711
$ws.Cells.Item(6, 6).Value = 2813
$ws.Cells.Item(7, 6).Value = 222
$ws.Cells.Item(8, 6).Value = 1550
$ws.Cells.Item(10, 6).Value = 6522
$ws.Cells.Item(11, 6).Value = 899
$ws.Cells.Item(12, 6).Value = 785
$ws.Cells.Item(13, 6).Value = 5083
$ws.Cells.Item(14, 6).Value = 4
$ws.Cells.Item(15, 6).Value = 563
$ws.Cells.Item(16, 6).Value = 2679
$ws.Cells.Item(17, 6).Value = 1368
$ws.Cells.Item(18, 6).Value = 1241
$ws.Cells.Item(21, 6).Value = 130
$ws.Cells.Item(22, 6).Value = 143
$ws.Cells.Item(23, 6).Value = 332
$ws.Cells.Item(24, 6).Value = 1110
$ws.Cells.Item(25, 6).Value = 259
$ws.Cells.Item(27, 6).Value = 159
$ws.Cells.Item(28, 6).Value = 553
$ws.Cells.Item(29, 6).Value = 1400
$ws.Cells.Item(31, 6).Value = 2122
$ws.Cells.Item(32, 6).Value = 599
$ws.Cells.Item(33, 6).Value = 49
$ws.Cells.Item(34, 6).Value = 166
$ws.Cells.Item(35, 6).Value = 44
$ws.Cells.Item(36, 6).Value = 267
$ws.Cells.Item(38, 6).Value = 1543
$ws.Cells.Item(40, 6).Value = 1101
$ws.Cells.Item(42, 6).Value = 31
$ws.Cells.Item(44, 6).Value = 2329
$ws.Cells.Item(45, 6).Value = 2593
$ws.Cells.Item(46, 6).Value = 152
